# address update test is added.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Log_in

# --- Add the new "My_address" worksheet right after "Log_in" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "My_address"

# ------------------------------------------------------------------
# Populate My_address (sheet2)
# ------------------------------------------------------------------

# Row 1 header (left part first)
$ws2.Range("A1").Value = "ID: 2"
$ws2.Range("B1").Value = "Test name: Update user address"

# Row 2
$ws2.Range("B2").Value = "ACTION"
$ws2.Range("C2").Value = "EXPECTED RESULT"
$ws2.Range("D2").Value = "DATA:"

# Row 3
$ws2.Range("B3").Value = "Pre conditions: user have an account"

# Row 4
$ws2.Range("B4").Value = "Go to ""http://automationpractice.com/index.php"""
$ws2.Range("C4").Value = "Web page is opened"
$ws2.Range("E4").Value = "Y"

# Row 5
$ws2.Range("B5").Value = "Click on ""sign in"" button on the upper right corner"
$ws2.Range("C5").Value = "Page witch sign in fields is opened"
$ws2.Range("E5").Value = "Y"

# Row 6
$ws2.Range("B6").Value = "Input email address in the email adress field"
$ws2.Range("C6").Value = "Email adress is inputed and visible"
$ws2.Range("D6").Value = "danka@fakemail.com"
$ws2.Hyperlinks.Add($ws2.Range("D6"), "mailto:danka@fakemail.com") | Out-Null
$ws2.Range("E6").Value = "Y"

# Row 7
$ws2.Range("B7").Value = "Input password in the password field"
$ws2.Range("C7").Value = "Password is inputed and visible"
$ws2.Range("D7").Value = 123456789
$ws2.Range("E7").Value = "Y"

# Row 8
$ws2.Range("B8").Value = "Click on ""sign in"" button "
$ws2.Range("C8").Value = "User is signed in"
$ws2.Range("E8").Value = "Y"

# --- Add "Bug report" column to the existing Log_in sheet ---
$ws1.Range("I1").Value = "Bug report"
$ws1.Range("I1").Font.Bold = $true
$ws1.Columns.Item(9).ColumnWidth = 9.6

# Finish row 1 header on My_address
$ws2.Range("E1").Value = "Test pass: Y/N"
$ws2.Range("F1").Value = "Comment"
$ws2.Range("G1").Value = "Bug report"

# Rows 9-10 (B & C columns)
$ws2.Range("B9").Value = "Click on MyAddresses button"
$ws2.Range("C9").Value = "My address page is opened"
$ws2.Range("B10").Value = "Click on ""update"" button"
$ws2.Range("C10").Value = "Address fields page is opened"

# Column E for rows 9-10
$ws2.Range("E9").Value = "y"
$ws2.Range("E10").Value = "Y"

# Row 11
$ws2.Range("B11").Value = "Input address2 in ""Address (Line 2)"" field"
$ws2.Range("D11").Value = "Novo naselje BB2"
$ws2.Range("C11").Value = "New address is inputed and visible"
$ws2.Range("E11").Value = "Y"

# Row 12
$ws2.Range("B12").Value = "Click on ""Save"" button"
$ws2.Range("D12").Value = "Your addresses are listed below."
$ws2.Range("C12").Value = "New address is saved"

Write-Host "done"
